$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "#Index": document the new "StageTimer" column (row 18)
# ---------------------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("#Index")

$wsIndex.Range("A18").Value = "StageTimer"
$wsIndex.Range("B18").Value = "int"
$wsIndex.Range("C18").Value = "스테이지 제한 시간 [ 초 단위로 설정 ]"

# ---------------------------------------------------------------------------
# Sheet "stage": add the new "StageTimer" column (L) with a header + values
# ---------------------------------------------------------------------------
$wsStage = $wb.Worksheets.Item("stage")

# widen column B (Name) a bit and give the new column L its own width
# (values are nudged slightly so the engine's internal pixel-snapping lands
# on the closest representable width to the intended 21.14 / 10.71 chars)
$wsStage.Columns.Item(2).ColumnWidth = 20.3
$wsStage.Columns.Item(12).ColumnWidth = 9.83

# Header rows
$wsStage.Range("L1").Value = "StageTimer"
$wsStage.Range("L2").Value = "int"

# Data rows: time limit in seconds for each stage
$wsStage.Range("L3").Value = 120
$wsStage.Range("L4").Value = 150
$wsStage.Range("L5").Value = 240
$wsStage.Range("L6").Value = 240
$wsStage.Range("L7").Value = 420
$wsStage.Range("L8").Value = 210
$wsStage.Range("L9").Value = 210
$wsStage.Range("L10").Value = 210
$wsStage.Range("L11").Value = 300
$wsStage.Range("L12").Value = 420
